$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 208, shifting existing rows 208..294 down to 209..295
$ws.Rows.Item(208).Insert()

# Populate the newly inserted row 208 with the new record's data
$ws.Cells.Item(208, 1).Value = 9
$ws.Cells.Item(208, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(208, 3).Value = "Metropolitana"
$ws.Cells.Item(208, 4).Value = 44466
$ws.Cells.Item(208, 5).Value = 13
$ws.Cells.Item(208, 6).Value = 100112013
$ws.Cells.Item(208, 7).Value = "Alcachofa"
$ws.Cells.Item(208, 8).Value = "Española"
$ws.Cells.Item(208, 9).Value = "Primera"
$ws.Cells.Item(208, 10).Value = 52
$ws.Cells.Item(208, 11).Value = 11000
$ws.Cells.Item(208, 12).Value = 12000
$ws.Cells.Item(208, 13).Value = 11500
$ws.Cells.Item(208, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(208, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(208, 16).Value = 383
$ws.Cells.Item(208, 17).Value = 30
$ws.Cells.Item(208, 18).Value = "Hortaliza"

# Ensure the date cell keeps the date number format used elsewhere in column D
$ws.Cells.Item(208, 4).NumberFormat = $ws.Cells.Item(209, 4).NumberFormat
